$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''27.900.96'
$ws.Range('E2').Value = '  +2.76%  '
$ws.Range('D3').Value = '''1.669.71'
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '''214.72'
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('E6').Value = '  -0.24%  '
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('D8').Value = '''23.63'
$ws.Range('E8').Value = '  +3.80%  '
$ws.Range('D9').Value = '''0.261'
$ws.Range('E9').Value = '  +0.21%  '
$ws.Range('E10').Value = '  +0.06%  '
$ws.Range('D11').Value = '''0.0878'
$ws.Range('E11').Value = '  -1.31%  '
$ws.Range('D12').Value = '''1.907.07'
$ws.Range('E12').Value = '  -0.28%  '
$ws.Range('D13').Value = '''1.664.01'
$ws.Range('E13').Value = '  -1.14%  '
$ws.Range('E14').Value = '  -1.29%  '
$ws.Range('E15').Value = '  -0.50%  '
$ws.Range('D16').Value = '''66.08'
$ws.Range('E16').Value = '  -0.76%  '
$ws.Range('D17').Value = '''251.74'
$ws.Range('E17').Value = '  +7.32%  '
$ws.Range('D18').Value = '''27.879.47'
$ws.Range('E18').Value = '  +2.90%  '
$ws.Range('D19').Value = '''0.0₃0731'
$ws.Range('E19').Value = '  -1.29%  '
$ws.Range('D20').Value = '''7.56'
$ws.Range('E20').Value = '  -3.91%  '
$ws.Range('D21').Value = '''0.999'
$ws.Range('E21').Value = '  -0.27%  '
$ws.Range('E22').Value = '  -1.53%  '
$ws.Range('D23').Value = '''9.36'
$ws.Range('E23').Value = '  -1.78%  '
$ws.Range('E24').Value = '  -1.60%  '
$ws.Range('D25').Value = '''146.68'
$ws.Range('E25').Value = '  -1.27%  '
$ws.Range('E26').Value = '  -2.95%  '
$ws.Range('D27').Value = '''16.31'
$ws.Range('E27').Value = '  -0.19%  '
$ws.Range('E28').Value = '  -0.18%  '
$ws.Range('E29').Value = '  -0.14%  '
$ws.Range('E30').Value = '  +6.29%  '
$ws.Range('D31').Value = '''0.0501'
$ws.Range('E31').Value = '  +0.41%  '
$ws.Range('E32').Value = '  -0.30%  '
$ws.Range('D33').Value = '''3.14'
$ws.Range('E33').Value = '  -2.49%  '
$ws.Range('D34').Value = '''1.424.25'
$ws.Range('E34').Value = '  -7.50%  '
$ws.Range('D35').Value = '''1.56'
$ws.Range('E35').Value = '  -5.36%  '
$ws.Range('D36').Value = '''2.39'
$ws.Range('E36').Value = '  +0.00%  '
$ws.Range('D37').Value = '''0.931'
$ws.Range('E37').Value = '  -1.05%  '
$ws.Range('D38').Value = '''0.583'
$ws.Range('E38').Value = '  -3.93%  '
$ws.Range('E39').Value = '  -0.96%  '
$ws.Range('B40').Value = 'WEMIXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D40').Value = '''1.03'
$ws.Range('E40').Value = '  -2.74%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value = '''69.63'
$ws.Range('E41').Value = '  +0.32%  '
$ws.Range('E42').Value = '  -0.24%  '
$ws.Range('E43').Value = '  -0.81%  '
$ws.Range('D44').Value = '''1.813.82'
$ws.Range('E44').Value = '  -0.88%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').Value = '''0.792'
$ws.Range('E45').Value = '  +1.63%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = '''5.38'
$ws.Range('E46').Value = '  -6.84%  '
$ws.Range('D47').Value = '''1.73'
$ws.Range('E47').Value = '  +5.48%  '
$ws.Range('D48').Value = '''88.89'
$ws.Range('E48').Value = '  -0.76%  '
$ws.Range('E49').Value = '  +0.75%  '
$ws.Range('D50').Value = '''0.102'
$ws.Range('E50').Value = '  -1.81%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '''7.79'
$ws.Range('E51').Value = '  -4.55%  '
